# Edit SSD.xlsx per diff:
#  - drop the duplicated "_x"/"_y" columns (I:M), keep a single column per metric
#  - rename/retarget remaining header cells (D1:H1)
#  - rewrite the data rows with the deduplicated + extended (2014-2021) dataset

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the redundant duplicate columns I:M entirely (shifts nothing afterwards left of them).
$ws.Range("I1:M1").EntireColumn.Delete()

# 2. Update the surviving header cells D1:H1 with the new consolidated labels.
$ws.Range("D1").Value = "4. Agriculture land area (% of land area)"
$ws.Range("E1").Value = "5. Average precipitation (mm per year)"
$ws.Range("F1").Value = "7. Fertilizer consumption (kilograms per hectare of arable land)"
$ws.Range("G1").Value = "13. Population"
$ws.Range("H1").Value = "17. Employment in agriculture (% of total employment) (modeled ILO estimate)"

# 3. Rewrite the data block (rows 2-9) with the new values (years 2014-2021).
$data = @(
    @("SSD", 2014, 98.25,  44.98494691, 900, 0.006753037,         11213284, 60.16503728435),
    @("SSD", 2015, 103.33, 44.92910607, 900, 0.024508935,         11194299, 60.3664853070102),
    @("SSD", 2016, 98.41,  44.87326524, 900, 0.034304838,         11066105, 60.9316367393641),
    @("SSD", 2017, 96.58,  44.8174244,  900, 0.026511029,         10658226, 61.1718858092903),
    @("SSD", 2018, 101.12, 44.70574272, 900, 0.10171903,          10395329, 61.2202485634641),
    @("SSD", 2019, 106.95, 44.70574272, 900, 0.089551927,         10447666, 60.9573437924524),
    @("SSD", 2020, 110.85, 44.70542623, 900, 0.07569632900000001, 10606227, 61.0425211864959),
    @("SSD", 2021, 112.4,  44.70859114, 900, 0.025844573,         10748272, 60.7832961686601)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $ws.Cells.Item($rowIndex, 8).Value = $row[7]
    $rowIndex = $rowIndex + 1
}

Write-Host "Final UsedRange:" $ws.UsedRange.Address()
